$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(537,1).NumberFormat = "@"
$ws.Cells.Item(537,1).Value = "2026/01/01"
$ws.Cells.Item(537,1).Style = "Normal"
$ws.Cells.Item(537,2).Value = "木"
$ws.Cells.Item(537,2).Style = "Normal"
$ws.Cells.Item(537,3).Value = 2
$ws.Cells.Item(537,3).Style = "Normal"
$ws.Cells.Item(537,4).Value = 13
$ws.Cells.Item(537,4).Style = "Normal"

$ws.Cells.Item(538,1).NumberFormat = "@"
$ws.Cells.Item(538,1).Value = "2026/01/01"
$ws.Cells.Item(538,1).Style = "Normal"
$ws.Cells.Item(538,2).Value = "木"
$ws.Cells.Item(538,2).Style = "Normal"
$ws.Cells.Item(538,3).Value = 5
$ws.Cells.Item(538,3).Style = "Normal"
$ws.Cells.Item(538,4).Value = 12
$ws.Cells.Item(538,4).Style = "Normal"

$ws.Cells.Item(539,1).NumberFormat = "@"
$ws.Cells.Item(539,1).Value = "2026/01/01"
$ws.Cells.Item(539,1).Style = "Normal"
$ws.Cells.Item(539,2).Value = "木"
$ws.Cells.Item(539,2).Style = "Normal"
$ws.Cells.Item(539,3).Value = 13
$ws.Cells.Item(539,3).Style = "Normal"
$ws.Cells.Item(539,4).Value = 14
$ws.Cells.Item(539,4).Style = "Normal"

$ws.Cells.Item(540,1).NumberFormat = "@"
$ws.Cells.Item(540,1).Value = "2026/01/01"
$ws.Cells.Item(540,1).Style = "Normal"
$ws.Cells.Item(540,2).Value = "木"
$ws.Cells.Item(540,2).Style = "Normal"
$ws.Cells.Item(540,3).Value = 16
$ws.Cells.Item(540,3).Style = "Normal"
$ws.Cells.Item(540,4).Value = 11
$ws.Cells.Item(540,4).Style = "Normal"

$ws.Cells.Item(541,1).NumberFormat = "@"
$ws.Cells.Item(541,1).Value = "2026/01/01"
$ws.Cells.Item(541,1).Style = "Normal"
$ws.Cells.Item(541,2).Value = "木"
$ws.Cells.Item(541,2).Style = "Normal"
$ws.Cells.Item(541,3).Value = 19
$ws.Cells.Item(541,3).Style = "Normal"
$ws.Cells.Item(541,4).Value = 13
$ws.Cells.Item(541,4).Style = "Normal"

$ws.Cells.Item(542,1).NumberFormat = "@"
$ws.Cells.Item(542,1).Value = "2026/01/02"
$ws.Cells.Item(542,1).Style = "Normal"
$ws.Cells.Item(542,2).Value = "金"
$ws.Cells.Item(542,2).Style = "Normal"
$ws.Cells.Item(542,3).Value = 1
$ws.Cells.Item(542,3).Style = "Normal"
$ws.Cells.Item(542,4).Value = 12
$ws.Cells.Item(542,4).Style = "Normal"

$ws.Cells.Item(543,1).NumberFormat = "@"
$ws.Cells.Item(543,1).Value = "2026/01/02"
$ws.Cells.Item(543,1).Style = "Normal"
$ws.Cells.Item(543,2).Value = "金"
$ws.Cells.Item(543,2).Style = "Normal"
$ws.Cells.Item(543,3).Value = 5
$ws.Cells.Item(543,3).Style = "Normal"
$ws.Cells.Item(543,4).Value = 12
$ws.Cells.Item(543,4).Style = "Normal"

$ws.Cells.Item(544,1).NumberFormat = "@"
$ws.Cells.Item(544,1).Value = "2026/01/02"
$ws.Cells.Item(544,1).Style = "Normal"
$ws.Cells.Item(544,2).Value = "金"
$ws.Cells.Item(544,2).Style = "Normal"
$ws.Cells.Item(544,3).Value = 8
$ws.Cells.Item(544,3).Style = "Normal"
$ws.Cells.Item(544,4).Value = 13
$ws.Cells.Item(544,4).Style = "Normal"

$ws.Cells.Item(545,1).NumberFormat = "@"
$ws.Cells.Item(545,1).Value = "2026/01/02"
$ws.Cells.Item(545,1).Style = "Normal"
$ws.Cells.Item(545,2).Value = "金"
$ws.Cells.Item(545,2).Style = "Normal"
$ws.Cells.Item(545,3).Value = 13
$ws.Cells.Item(545,3).Style = "Normal"
$ws.Cells.Item(545,4).Value = 16
$ws.Cells.Item(545,4).Style = "Normal"

$ws.Cells.Item(546,1).NumberFormat = "@"
$ws.Cells.Item(546,1).Value = "2026/01/02"
$ws.Cells.Item(546,1).Style = "Normal"
$ws.Cells.Item(546,2).Value = "金"
$ws.Cells.Item(546,2).Style = "Normal"
$ws.Cells.Item(546,3).Value = 16
$ws.Cells.Item(546,3).Style = "Normal"
$ws.Cells.Item(546,4).Value = 19
$ws.Cells.Item(546,4).Style = "Normal"

$ws.Cells.Item(547,1).NumberFormat = "@"
$ws.Cells.Item(547,1).Value = "2026/01/02"
$ws.Cells.Item(547,1).Style = "Normal"
$ws.Cells.Item(547,2).Value = "金"
$ws.Cells.Item(547,2).Style = "Normal"
$ws.Cells.Item(547,3).Value = 19
$ws.Cells.Item(547,3).Style = "Normal"
$ws.Cells.Item(547,4).Value = 21
$ws.Cells.Item(547,4).Style = "Normal"

$ws.Cells.Item(548,1).NumberFormat = "@"
$ws.Cells.Item(548,1).Value = "2026/01/02"
$ws.Cells.Item(548,1).Style = "Normal"
$ws.Cells.Item(548,2).Value = "金"
$ws.Cells.Item(548,2).Style = "Normal"
$ws.Cells.Item(548,3).Value = 22
$ws.Cells.Item(548,3).Style = "Normal"
$ws.Cells.Item(548,4).Value = 22
$ws.Cells.Item(548,4).Style = "Normal"

$ws.Cells.Item(549,1).NumberFormat = "@"
$ws.Cells.Item(549,1).Value = "2026/01/03"
$ws.Cells.Item(549,1).Style = "Normal"
$ws.Cells.Item(549,2).Value = "土"
$ws.Cells.Item(549,2).Style = "Normal"
$ws.Cells.Item(549,3).Value = 1
$ws.Cells.Item(549,3).Style = "Normal"
$ws.Cells.Item(549,4).Value = 23
$ws.Cells.Item(549,4).Style = "Normal"

$ws.Cells.Item(550,1).NumberFormat = "@"
$ws.Cells.Item(550,1).Value = "2026/01/03"
$ws.Cells.Item(550,1).Style = "Normal"
$ws.Cells.Item(550,2).Value = "土"
$ws.Cells.Item(550,2).Style = "Normal"
$ws.Cells.Item(550,3).Value = 4
$ws.Cells.Item(550,3).Style = "Normal"
$ws.Cells.Item(550,4).Value = 26
$ws.Cells.Item(550,4).Style = "Normal"

$ws.Cells.Item(551,1).NumberFormat = "@"
$ws.Cells.Item(551,1).Value = "2026/01/03"
$ws.Cells.Item(551,1).Style = "Normal"
$ws.Cells.Item(551,2).Value = "土"
$ws.Cells.Item(551,2).Style = "Normal"
$ws.Cells.Item(551,3).Value = 7
$ws.Cells.Item(551,3).Style = "Normal"
$ws.Cells.Item(551,4).Value = 23
$ws.Cells.Item(551,4).Style = "Normal"

$ws.Cells.Item(552,1).NumberFormat = "@"
$ws.Cells.Item(552,1).Value = "2026/01/03"
$ws.Cells.Item(552,1).Style = "Normal"
$ws.Cells.Item(552,2).Value = "土"
$ws.Cells.Item(552,2).Style = "Normal"
$ws.Cells.Item(552,3).Value = 13
$ws.Cells.Item(552,3).Style = "Normal"
$ws.Cells.Item(552,4).Value = 23
$ws.Cells.Item(552,4).Style = "Normal"

$ws.Cells.Item(553,1).NumberFormat = "@"
$ws.Cells.Item(553,1).Value = "2026/01/03"
$ws.Cells.Item(553,1).Style = "Normal"
$ws.Cells.Item(553,2).Value = "土"
$ws.Cells.Item(553,2).Style = "Normal"
$ws.Cells.Item(553,3).Value = 16
$ws.Cells.Item(553,3).Style = "Normal"
$ws.Cells.Item(553,4).Value = 24
$ws.Cells.Item(553,4).Style = "Normal"

$ws.Cells.Item(554,1).NumberFormat = "@"
$ws.Cells.Item(554,1).Value = "2026/01/03"
$ws.Cells.Item(554,1).Style = "Normal"
$ws.Cells.Item(554,2).Value = "土"
$ws.Cells.Item(554,2).Style = "Normal"
$ws.Cells.Item(554,3).Value = 19
$ws.Cells.Item(554,3).Style = "Normal"
$ws.Cells.Item(554,4).Value = 26
$ws.Cells.Item(554,4).Style = "Normal"

$ws.Cells.Item(555,1).NumberFormat = "@"
$ws.Cells.Item(555,1).Value = "2026/01/03"
$ws.Cells.Item(555,1).Style = "Normal"
$ws.Cells.Item(555,2).Value = "土"
$ws.Cells.Item(555,2).Style = "Normal"
$ws.Cells.Item(555,3).Value = 22
$ws.Cells.Item(555,3).Style = "Normal"
$ws.Cells.Item(555,4).Value = 21
$ws.Cells.Item(555,4).Style = "Normal"

$ws.Cells.Item(556,1).NumberFormat = "@"
$ws.Cells.Item(556,1).Value = "2026/01/04"
$ws.Cells.Item(556,1).Style = "Normal"
$ws.Cells.Item(556,2).Value = "日"
$ws.Cells.Item(556,2).Style = "Normal"
$ws.Cells.Item(556,3).Value = 2
$ws.Cells.Item(556,3).Style = "Normal"
$ws.Cells.Item(556,4).Value = 19
$ws.Cells.Item(556,4).Style = "Normal"

$ws.Cells.Item(557,1).NumberFormat = "@"
$ws.Cells.Item(557,1).Value = "2026/01/04"
$ws.Cells.Item(557,1).Style = "Normal"
$ws.Cells.Item(557,2).Value = "日"
$ws.Cells.Item(557,2).Style = "Normal"
$ws.Cells.Item(557,3).Value = 4
$ws.Cells.Item(557,3).Style = "Normal"
$ws.Cells.Item(557,4).Value = 18
$ws.Cells.Item(557,4).Style = "Normal"

$ws.Cells.Item(558,1).NumberFormat = "@"
$ws.Cells.Item(558,1).Value = "2026/01/04"
$ws.Cells.Item(558,1).Style = "Normal"
$ws.Cells.Item(558,2).Value = "日"
$ws.Cells.Item(558,2).Style = "Normal"
$ws.Cells.Item(558,3).Value = 7
$ws.Cells.Item(558,3).Style = "Normal"
$ws.Cells.Item(558,4).Value = 19
$ws.Cells.Item(558,4).Style = "Normal"

$ws.Cells.Item(559,1).NumberFormat = "@"
$ws.Cells.Item(559,1).Value = "2026/01/04"
$ws.Cells.Item(559,1).Style = "Normal"
$ws.Cells.Item(559,2).Value = "日"
$ws.Cells.Item(559,2).Style = "Normal"
$ws.Cells.Item(559,3).Value = 13
$ws.Cells.Item(559,3).Style = "Normal"
$ws.Cells.Item(559,4).Value = 20
$ws.Cells.Item(559,4).Style = "Normal"

$ws.Cells.Item(560,1).NumberFormat = "@"
$ws.Cells.Item(560,1).Value = "2026/01/04"
$ws.Cells.Item(560,1).Style = "Normal"
$ws.Cells.Item(560,2).Value = "日"
$ws.Cells.Item(560,2).Style = "Normal"
$ws.Cells.Item(560,3).Value = 22
$ws.Cells.Item(560,3).Style = "Normal"
$ws.Cells.Item(560,4).Value = 13
$ws.Cells.Item(560,4).Style = "Normal"

$ws.Cells.Item(561,1).NumberFormat = "@"
$ws.Cells.Item(561,1).Value = "2026/01/05"
$ws.Cells.Item(561,1).Style = "Normal"
$ws.Cells.Item(561,2).Value = "月"
$ws.Cells.Item(561,2).Style = "Normal"
$ws.Cells.Item(561,3).Value = 1
$ws.Cells.Item(561,3).Style = "Normal"
$ws.Cells.Item(561,4).Value = 13
$ws.Cells.Item(561,4).Style = "Normal"

$ws.Cells.Item(562,1).NumberFormat = "@"
$ws.Cells.Item(562,1).Value = "2026/01/05"
$ws.Cells.Item(562,1).Style = "Normal"
$ws.Cells.Item(562,2).Value = "月"
$ws.Cells.Item(562,2).Style = "Normal"
$ws.Cells.Item(562,3).Value = 7
$ws.Cells.Item(562,3).Style = "Normal"
$ws.Cells.Item(562,4).Value = 14
$ws.Cells.Item(562,4).Style = "Normal"
